# Update F-column ('想去人数' / want-to-go headcount) values on the
# '展览' sheet (index 1) and '全部类型' sheet (index 4), per the
# upstream data refresh (gh-pages output regenerated at 456a3b4).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws1.Cells.Item(3, 6).Value = 581   # F3: 579 -> 581
$ws1.Cells.Item(5, 6).Value = 295   # F5: 294 -> 295
$ws1.Cells.Item(6, 6).Value = 1115   # F6: 1113 -> 1115
$ws1.Cells.Item(7, 6).Value = 1449   # F7: 1448 -> 1449
$ws1.Cells.Item(9, 6).Value = 116   # F9: 115 -> 116
$ws1.Cells.Item(10, 6).Value = 757   # F10: 754 -> 757
$ws1.Cells.Item(13, 6).Value = 120   # F13: 121 -> 120
$ws1.Cells.Item(14, 6).Value = 451   # F14: 450 -> 451
$ws1.Cells.Item(15, 6).Value = 1392   # F15: 1391 -> 1392
$ws1.Cells.Item(16, 6).Value = 125   # F16: 124 -> 125
$ws1.Cells.Item(18, 6).Value = 282   # F18: 281 -> 282
$ws1.Cells.Item(20, 6).Value = 75   # F20: 74 -> 75
$ws1.Cells.Item(21, 6).Value = 663   # F21: 662 -> 663
$ws1.Cells.Item(22, 6).Value = 1014   # F22: 1013 -> 1014
$ws1.Cells.Item(24, 6).Value = 247   # F24: 243 -> 247
$ws1.Cells.Item(26, 6).Value = 6003   # F26: 5988 -> 6003
$ws1.Cells.Item(31, 6).Value = 14703   # F31: 14690 -> 14703
$ws1.Cells.Item(32, 6).Value = 1459   # F32: 1457 -> 1459
$ws1.Cells.Item(33, 6).Value = 231   # F33: 228 -> 231
$ws1.Cells.Item(36, 6).Value = 9497   # F36: 9488 -> 9497
$ws1.Cells.Item(37, 6).Value = 650   # F37: 647 -> 650
$ws1.Cells.Item(38, 6).Value = 4227   # F38: 4225 -> 4227
$ws1.Cells.Item(39, 6).Value = 162   # F39: 161 -> 162

$ws4 = $wb.Worksheets.Item(4)   # 全部类型
$ws4.Cells.Item(3, 6).Value = 581   # F3: 579 -> 581
$ws4.Cells.Item(5, 6).Value = 295   # F5: 294 -> 295
$ws4.Cells.Item(6, 6).Value = 1115   # F6: 1113 -> 1115
$ws4.Cells.Item(7, 6).Value = 1449   # F7: 1448 -> 1449
$ws4.Cells.Item(9, 6).Value = 116   # F9: 115 -> 116
$ws4.Cells.Item(10, 6).Value = 757   # F10: 754 -> 757
$ws4.Cells.Item(13, 6).Value = 120   # F13: 121 -> 120
$ws4.Cells.Item(14, 6).Value = 451   # F14: 450 -> 451
$ws4.Cells.Item(15, 6).Value = 1392   # F15: 1391 -> 1392
$ws4.Cells.Item(16, 6).Value = 125   # F16: 124 -> 125
$ws4.Cells.Item(18, 6).Value = 282   # F18: 281 -> 282
$ws4.Cells.Item(21, 6).Value = 75   # F21: 74 -> 75
$ws4.Cells.Item(22, 6).Value = 663   # F22: 662 -> 663
$ws4.Cells.Item(24, 6).Value = 1014   # F24: 1013 -> 1014
$ws4.Cells.Item(26, 6).Value = 247   # F26: 243 -> 247
$ws4.Cells.Item(29, 6).Value = 6003   # F29: 5988 -> 6003
$ws4.Cells.Item(34, 6).Value = 14703   # F34: 14690 -> 14703
$ws4.Cells.Item(35, 6).Value = 1459   # F35: 1457 -> 1459
$ws4.Cells.Item(36, 6).Value = 231   # F36: 228 -> 231
$ws4.Cells.Item(39, 6).Value = 9497   # F39: 9488 -> 9497
$ws4.Cells.Item(40, 6).Value = 650   # F40: 647 -> 650
$ws4.Cells.Item(41, 6).Value = 4227   # F41: 4225 -> 4227
$ws4.Cells.Item(42, 6).Value = 162   # F42: 161 -> 162
